# Rename 'Codelists' sheet to 'Cells', and make it the active/selected sheet
# with its selection set to F13 (previously the 'Variables' sheet was the
# selected tab, and the 'Codelists' sheet's own selection was left at H27).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Codelists")
$ws.Name = "Cells"

# Make the renamed sheet the active tab (this also clears tabSelected on
# whichever sheet was previously active, i.e. "Variables").
$ws.Activate() | Out-Null

# Update the selection on the renamed sheet.
$ws.Range("F13").Select() | Out-Null
